# Update Excel file from GitHub Actions on 2025-10-25 15:17:01
#
# Appends the next day's gold-price row (25-10-2025) to the bottom of the
# single data table on Sheet1, mirroring the formatting of the row above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the first empty row right after the existing data (row 40 -> 41).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Carry over the formatting (borders, wrap-text, etc.) of the row above,
# exactly like the previous day's rows were extended.
$ws.Range("A" + $lastRow + ":B" + $lastRow).Copy()
$ws.Range("A" + $newRow + ":B" + $newRow).PasteSpecial(-4122)

# Write the new day's values.
$ws.Cells.Item($newRow, 1).Value = "25-10-2025"
$ws.Cells.Item($newRow, 2).Value = "The price of gold in India today is ₹12,562 per gram for 24 karat gold, ₹11,515 per gram for 22 karat gold and ₹9,422 per gram for 18 karat gold (also called 999 gold)."
